# The workbook records one "expense/transaction" entry per row on the
# "2024" sheet. A new September entry was logged ("balance your axis" at
# 2024-09-03 13:14:06), which was inserted as a brand-new row right above
# the existing "lounge" / 2024-09-03 13:08:08 row (row 29), pushing every
# row below it (29-40) down by one (to 30-41) and extending the used
# range from A1:Y40 to A1:Y41.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row above the current row 29, shifting rows 29:40
# down to 30:41 (this also grows the sheet's dimension to A1:Y41).
$ws.Rows("29").Insert()

# Populate the new row 29 with the new September transaction.
$ws.Range("R29").Value = "balance your axis"
$ws.Range("S29").Value = "2024-09-03 13:14:06"
